# Update the frequency table (rows 2-5, columns B-X) with the re-run
# publication values. Only the cells whose value actually changed are
# written here; unchanged cells (several zeros) are left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.239322533136966
$ws.Range("C2").Value = 0.0331369661266568
$ws.Range("D2").Value = 0.0147275405007364
$ws.Range("E2").Value = 0.00368188512518409
$ws.Range("F2").Value = 0.166421207658321
$ws.Range("G2").Value = 0.00294550810014728
$ws.Range("H2").Value = 0.000736377025036819
$ws.Range("I2").Value = 0.00662739322533137
$ws.Range("J2").Value = 0.0176730486008837
$ws.Range("K2").Value = 0.00441826215022091
$ws.Range("L2").Value = 0.00957290132547865
$ws.Range("M2").Value = 0.846833578792342
$ws.Range("N2").Value = 0.00515463917525773
$ws.Range("O2").Value = 0.00957290132547865
$ws.Range("P2").Value = 0.0191458026509573
$ws.Range("Q2").Value = 0.00368188512518409
$ws.Range("S2").Value = 0.00441826215022091
$ws.Range("T2").Value = 0.993372606774669
$ws.Range("U2").Value = 0.0103092783505155
$ws.Range("V2").Value = 0.00368188512518409
$ws.Range("W2").Value = 0.00957290132547865
$ws.Range("X2").Value = 0.00220913107511046
$ws.Range("B3").Value = 0.0184094256259205
$ws.Range("C3").Value = 0.770986745213549
$ws.Range("D3").Value = 0.818851251840943
$ws.Range("E3").Value = 0.00441826215022091
$ws.Range("F3").Value = 0.00736377025036819
$ws.Range("G3").Value = 0.164948453608247
$ws.Range("H3").Value = 0.00515463917525773
$ws.Range("I3").Value = 0.815905743740795
$ws.Range("J3").Value = 0.00736377025036819
$ws.Range("K3").Value = 0.156848306332842
$ws.Range("L3").Value = 0.156111929307806
$ws.Range("M3").Value = 0.00147275405007364
$ws.Range("N3").Value = 0.981590574374079
$ws.Range("O3").Value = 0.00441826215022091
$ws.Range("P3").Value = 0.000736377025036819
$ws.Range("Q3").Value = 0.00220913107511046
$ws.Range("R3").Value = 0.00810014727540501
$ws.Range("U3").Value = 0.1620029455081
$ws.Range("V3").Value = 0.0125184094256259
$ws.Range("W3").Value = 0.0147275405007364
$ws.Range("X3").Value = 0.0139911634756996
$ws.Range("B4").Value = 0.733431516936672
$ws.Range("C4").Value = 0.163475699558174
$ws.Range("D4").Value = 0.1620029455081
$ws.Range("E4").Value = 0.181885125184094
$ws.Range("F4").Value = 0.823269513991163
$ws.Range("G4").Value = 0.822533136966127
$ws.Range("H4").Value = 0.835051546391753
$ws.Range("I4").Value = 0.0103092783505155
$ws.Range("J4").Value = 0.815905743740795
$ws.Range("K4").Value = 0.831369661266569
$ws.Range("L4").Value = 0.832106038291605
$ws.Range("M4").Value = 0.150220913107511
$ws.Range("N4").Value = 0.00589101620029455
$ws.Range("O4").Value = 0.000736377025036819
$ws.Range("P4").Value = 0.980117820324006
$ws.Range("Q4").Value = 0.993372606774669
$ws.Range("R4").Value = 0.988954344624448
$ws.Range("S4").Value = 0.995581737849779
$ws.Range("T4").Value = 0.00368188512518409
$ws.Range("U4").Value = 0.826951399116348
$ws.Range("V4").Value = 0.0235640648011782
$ws.Range("W4").Value = 0.814432989690722
$ws.Range("X4").Value = 0.824742268041237
$ws.Range("B5").Value = 0.00810014727540501
$ws.Range("C5").Value = 0.0309278350515464
$ws.Range("D5").Value = 0.00441826215022091
$ws.Range("E5").Value = 0.810014727540501
$ws.Range("F5").Value = 0.00294550810014728
$ws.Range("G5").Value = 0.00957290132547865
$ws.Range("H5").Value = 0.157584683357879
$ws.Range("I5").Value = 0.167157584683358
$ws.Range("J5").Value = 0.159057437407953
$ws.Range("K5").Value = 0.00736377025036819
$ws.Range("L5").Value = 0.00220913107511046
$ws.Range("M5").Value = 0.00147275405007364
$ws.Range("N5").Value = 0.00736377025036819
$ws.Range("O5").Value = 0.985272459499264
$ws.Range("Q5").Value = 0.000736377025036819
$ws.Range("R5").Value = 0.00294550810014728
$ws.Range("T5").Value = 0.00294550810014728
$ws.Range("U5").Value = 0.000736377025036819
$ws.Range("V5").Value = 0.960235640648012
$ws.Range("W5").Value = 0.161266568483063
$ws.Range("X5").Value = 0.158321060382916

Write-Output "done"
